$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 2121
$ws1.Range("F4").Value = 873
$ws1.Range("F5").Value = 1343
$ws1.Range("F6").Value = 366

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 2121
$ws4.Range("F6").Value = 873
$ws4.Range("F7").Value = 1343
$ws4.Range("F8").Value = 366
